$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the column E width to match column D
$ws.Columns.Item(5).ColumnWidth = 12.83203125

# Add the new column E data (force text so date-like strings aren't
# auto-converted to date serial values, matching column D's behaviour)
$ws.Range("E1:E3").NumberFormat = "@"
$ws.Range("E1").Value = "02/01/2555"
$ws.Range("E2").Value = "ขาดเรียน"
$ws.Range("E3").Value = "มาเรียน"
